$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.845.50'
$ws.Range("E2").Value = '  +2.42%  '
$ws.Range("D3").Value = '1.863.84'
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '247.05'
$ws.Range("E5").Value = '  +2.15%  '
$ws.Range("D6").Value = '0.6386'
$ws.Range("E6").Value = '  +3.42%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '0.3007'
$ws.Range("E8").Value = '  +4.10%  '
$ws.Range("D9").Value = '0.07503'
$ws.Range("E9").Value = '  +2.07%  '
$ws.Range("D10").Value = '24.26'
$ws.Range("E10").Value = '  +5.71%  '
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").Value = '1.897.62'
$ws.Range("E12").Value = '  +4.23%  '
$ws.Range("D13").Value = '5.073'
$ws.Range("E13").Value = '  +2.27%  '
$ws.Range("D14").Value = '0.6916'
$ws.Range("E14").Value = '  +4.44%  '
$ws.Range("D15").Value = '84.72'
$ws.Range("E15").Value = '  +3.61%  '
$ws.Range("D16").Value = '0.000009512'
$ws.Range("E16").Value = '  +6.23%  '
$ws.Range("D17").Value = '6.137'
$ws.Range("E17").Value = '  +4.73%  '
$ws.Range("D18").Value = '29.843.82'
$ws.Range("E18").Value = '  +2.50%  '
$ws.Range("D19").Value = '2.103.67'
$ws.Range("E19").Value = '  +2.55%  '
$ws.Range("D20").Value = '241.17'
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("E21").Value = '  +1.81%  '
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = '7.396'
$ws.Range("E23").Value = '  +3.54%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '159.57'
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").Value = '0.1430'
$ws.Range("E26").Value = '  +1.47%  '
$ws.Range("D27").Value = '8.587'
$ws.Range("E27").Value = '  +1.60%  '
$ws.Range("D28").Value = '18.02'
$ws.Range("E28").Value = '  +2.08%  '
$ws.Range("D29").Value = '1.510'
$ws.Range("E29").Value = '  +1.87%  '
$ws.Range("D30").Value = '0.06022'
$ws.Range("E30").Value = '  +7.70%  '
$ws.Range("D31").Value = '1.266'
$ws.Range("E31").Value = '  +4.95%  '
$ws.Range("D32").Value = '4.147'
$ws.Range("E32").Value = '  +1.22%  '
$ws.Range("D33").Value = '4.153'
$ws.Range("E33").Value = '  +1.09%  '
$ws.Range("D34").Value = '1.879'
$ws.Range("E34").Value = '  +2.97%  '
$ws.Range("E35").Value = '  +2.65%  '
$ws.Range("D36").Value = '0.7370'
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("D37").Value = '2.615'
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("D38").Value = '2.879'
$ws.Range("E38").Value = '  +1.59%  '
$ws.Range("D39").Value = '1.228.25'
$ws.Range("E39").Value = '  +1.69%  '
$ws.Range("D40").Value = '0.01796'
$ws.Range("E40").Value = '  +1.92%  '
$ws.Range("D41").Value = '6.410'
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("D42").Value = '0.9255'
$ws.Range("E42").Value = '  +3.71%  '
$ws.Range("D43").Value = '1.002'
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '102.69'
$ws.Range("E44").Value = '  +1.98%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '2.004.92'
$ws.Range("E45").Value = '  +2.52%  '
$ws.Range("D46").Value = '66.70'
$ws.Range("E46").Value = '  +3.00%  '
$ws.Range("D47").Value = '0.00000000121'
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").Value = '0.5084'
$ws.Range("D49").Value = '9.346'
$ws.Range("E49").Value = '  +2.90%  '
$ws.Range("D50").Value = '0.4105'
$ws.Range("E50").Value = '  +2.67%  '
$ws.Range("D51").Value = '0.1147'
$ws.Range("E51").Value = '  +3.21%  '
